# Card1: fix the "Correction " header (drop the trailing space) and add
# a new "Serviced by " column (O) after it, following the same layout
# already used on the sibling Card sheets (Card15/Card22/Card7/Card2).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# N1 header: "Correction " -> "Correction"
$ws.Cells.Item(1, 14).Value = "Correction"

# New O1 header, cloned from N1's formatting (bold / border / centered)
# so it matches the rest of row 1.
$ws.Cells.Item(1, 14).Copy() | Out-Null
$ws.Cells.Item(1, 15).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 15).Value = "Serviced by "

# Data rows 2-12: the N column already held "nan" placeholders in every
# other column, so restore that here; O starts out blank (text-typed,
# like the rest of the sheet's placeholder cells) with no special
# formatting applied.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Cells.Item($r, 15).Style = "Normal"
}
